$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values for the "Price" column (D)
# that look like plain numbers are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the original inlineStr/text cells)
# instead of silently converting them to numeric values.
$updates = @{
    "D2" = "64.441.48"; "E2" = "  -2.81%  "
    "D3" = "3.175.57"; "E3" = "  -4.44%  "
    "E4" = "  +0.03%  "
    "D5" = "'570.85"; "E5" = "  -2.81%  "
    "D6" = "'168.76"; "E6" = "  -8.30%  "
    "D7" = "'0.606"; "E7" = "  -6.63%  "
    "E8" = "  -0.11%  "
    "D9" = "3.186.61"; "E9" = "  -4.09%  "
    "E10" = "  -4.08%  "
    "D11" = "'6.81"; "E11" = "  -0.30%  "
    "E12" = "  -3.66%  "
    "D13" = "3.740.72"; "E13" = "  -4.12%  "
    "E14" = "  -2.17%  "
    "D15" = "64.501.04"; "E15" = "  -2.76%  "
    "D16" = "'25.39"; "E16" = "  -3.08%  "
    "E17" = "  -2.99%  "
    "D18" = "3.186.35"; "E18" = "  -4.20%  "
    "D19" = "'415.84"; "E19" = "  -2.23%  "
    "E20" = "  -2.10%  "
    "E21" = "  -3.48%  "
    "E22" = "  -3.95%  "
    "E23" = "  -0.26%  "
    "D24" = "'70.47"; "E24" = "  -1.93%  "
    "D25" = "'5.67"; "E25" = "  +0.02%  "
    "E26" = "  +0.92%  "
    "D27" = "'0.488"; "E27" = "  -5.13%  "
    "E28" = "  -7.00%  "
    "D29" = "'8.77"; "E29" = "  -1.81%  "
    "D30" = "'0.999"; "E30" = "  -0.08%  "
    "E31" = "  -3.95%  "
    "E32" = "  -2.92%  "
    "D34" = "'5.10"; "E34" = "  -1.69%  "
    "E35" = "  -3.83%  "
    "E36" = "  -3.86%  "
    "D37" = "'158.55"; "E37" = "  -0.79%  "
    "E38" = "  -5.47%  "
    "D39" = "2.733.09"; "E39" = "  -5.45%  "
    "E40" = "  -5.51%  "
    "D41" = "'24.46"; "E41" = "  -7.50%  "
    "E42" = "  -2.80%  "
    "E43" = "  -2.19%  "
    "E44" = "  -6.44%  "
    "E45" = "  -6.27%  "
    "D46" = "'5.62"; "E46" = "  -5.60%  "
    "E47" = "  -2.89%  "
    "D48" = "'21.69"; "E48" = "  -6.77%  "
    "D49" = "'293.37"; "E49" = "  -6.79%  "
    "D50" = "'2.01"; "E50" = "  -12.70%  "
    "E51" = "  -0.18%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
